$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- USUARIOS_TIPOS table (was USUARIO_TIPO): title + header ---
$ws.Range("H7").Value = "USUARIOS_TIPOS"
$ws.Range("I8").Value = "Tipo"

# --- Small table header swap near STATUS_CONSULTAS / GÊNEROS ---
$ws.Range("C26").Value = "Tipo"
$ws.Range("F26").Value = "Nome"

# --- PRONTUARIOS table: Gênero -> Id_Gênero (must land in shared-strings
#     table before "médico" below) ---
$ws.Range("B40").Value = "Id_Gênero"

# --- USUARIOS_TIPOS data row ---
$ws.Range("I10").Value = "médico"

# --- Update view: scroll so F19 is top-left, select I10 ---
$ws.Range("I10").Select()
